# TOLLCLASS_Designations.xlsx - restore changes that were overwritten
# (commit: "changes to TOLLCLASS_Designations.xlsx in a previous commit got overwritten")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs_for_tollcalib")

# 1) Header rename: "PBA2050 RTP ID" -> "PBA2050_RTP_ID"
$ws.Range("D1").Value = "PBA2050_RTP_ID"

# 2) Populate toll-setting defaults (F=Toll45, G=5, H=0.03) for the newly
#    added NextGenFwyR2 rows 466:533 that previously only had A:C filled in.
$ws.Range("F466:F533").Value = 45
$ws.Range("G466:G533").Value = 5
$ws.Range("H466:H533").Value = 0.03

# 3) Restore the view/selection state that was current when the author saved.
[void]$ws.Range("D2").Select()
$excel.ActiveWindow.ScrollRow = 2
